# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 344
$ws.Range("F3").Value  = 238
$ws.Range("F4").Value  = 560
$ws.Range("F5").Value  = 1347
$ws.Range("F6").Value  = 655
$ws.Range("F7").Value  = 350
$ws.Range("F8").Value  = 31
$ws.Range("F10").Value = 413
$ws.Range("F11").Value = 6219
$ws.Range("F15").Value = 4649
$ws.Range("F16").Value = 461
$ws.Range("F19").Value = 5428
$ws.Range("F20").Value = 7098
$ws.Range("F21").Value = 148
$ws.Range("F22").Value = 1087
$ws.Range("F23").Value = 752
$ws.Range("F24").Value = 3992
$ws.Range("F25").Value = 552
$ws.Range("F28").Value = 145
$ws.Range("F29").Value = 1057
$ws.Range("F30").Value = 1492
$ws.Range("F31").Value = 552
$ws.Range("F32").Value = 684
$ws.Range("F33").Value = 1685
$ws.Range("F35").Value = 1879
$ws.Range("F36").Value = 229
$ws.Range("F37").Value = 41
$ws.Range("F38").Value = 1239
$ws.Range("F40").Value = 685
$ws.Range("F41").Value = 319
$ws.Range("F42").Value = 1009
$ws.Range("F43").Value = 3659
$ws.Range("F45").Value = 342
$ws.Range("F46").Value = 443
$ws.Range("F47").Value = 23
$ws.Range("F48").Value = 91
$ws.Range("F49").Value = 3951

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 1265
$ws.Range("F28").Value = 86
$ws.Range("F29").Value = 24

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4387

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 4387
$ws.Range("F3").Value  = 344
$ws.Range("F4").Value  = 1265
$ws.Range("F6").Value  = 238
$ws.Range("F7").Value  = 560
$ws.Range("F9").Value  = 1347
$ws.Range("F11").Value = 655
$ws.Range("F12").Value = 350
$ws.Range("F14").Value = 413
$ws.Range("F17").Value = 4649
$ws.Range("F18").Value = 5428
$ws.Range("F19").Value = 5428
$ws.Range("F20").Value = 148
$ws.Range("F21").Value = 1087
$ws.Range("F22").Value = 752
$ws.Range("F23").Value = 3992
$ws.Range("F24").Value = 552
$ws.Range("F27").Value = 145
$ws.Range("F28").Value = 1057
$ws.Range("F29").Value = 1492
$ws.Range("F30").Value = 552
$ws.Range("F31").Value = 684
$ws.Range("F32").Value = 1685
$ws.Range("F34").Value = 1879
$ws.Range("F39").Value = 685
$ws.Range("F41").Value = 319
$ws.Range("F42").Value = 86
$ws.Range("F43").Value = 3659
$ws.Range("F44").Value = 24
$ws.Range("F46").Value = 342
$ws.Range("F47").Value = 443
$ws.Range("F48").Value = 91
$ws.Range("F50").Value = 3951
